$wb = $excel.ActiveWorkbook

# Add the new "Calls" worksheet after the last existing sheet (Cases)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$newSheet.Name = "Calls"

# Header row fill/highlight (same yellow highlight style used on other sheets' headers)
$newSheet.Range("A1:E1").Interior.Color = 65535

# Header row (left to right)
$newSheet.Range("A1").Value = "contact"
$newSheet.Range("B1").Value = "deal"
$newSheet.Range("C1").Value = "task"
$newSheet.Range("D1").Value = "case"
$newSheet.Range("E1").Value = "notes"

# Data filled column by column
$newSheet.Range("A2").Value = "zzzx"
$newSheet.Range("A3").Value = "aaaa"

$newSheet.Range("B2").Value = "aaaa"
$newSheet.Range("B3").Value = "bbbb"

$newSheet.Range("C2").Value = "cccc"
$newSheet.Range("C3").Value = "dddd"

$newSheet.Range("D2").Value = "eeee"
$newSheet.Range("D3").Value = "ffff"

$newSheet.Range("E2").Value = "gggg"
$newSheet.Range("E3").Value = "hhhh"

# The previously active sheet (Cases) no longer has a single-cell selection
$wb.Worksheets.Item("Cases").Range("A1:XFD1").Select() | Out-Null

# Make the new Calls sheet the active/selected sheet, with column B selected
$newSheet.Activate() | Out-Null
$newSheet.Range("B1:B1048576").Select() | Out-Null
